$d = $word.ActiveDocument

# --- Edit 1: Merge "Also" sentence runs into a single run (remove proofErr split) ---
$rng1 = $d.Content
$rng1.Find.Execute("Previous reviews of literature", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para1 = $rng1.Paragraphs(1).Range
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Previous reviews of literature that are extensive in nature would have provided a better theoretical basis and hypothesis development. There were some early technical issues related to dependencies of the R packages that led to unforeseen delays in the analysis. The score of the evidence outlier 5 could have been detected earlier in the exploration analysis, which can guide the initial hypotheses. Also having investigated possible confounding factors such as intensity of media coverage, prevalence of the condition and celebrity endorsement would have given more mechanistic understanding of the weak correlation despite the big differences between the groups.</w:t></w:r></w:p>'
$para1.InsertXML($xml1)

# --- Edit 2: Split "-Wallis" run into 3 runs with proofErr gramStart/gramEnd around "Wallis" ---
$rng2 = $d.Content
$rng2.Find.Execute("The GitHub repository shows systematic development", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para2 = $rng2.Paragraphs(1).Range
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">The GitHub repository shows systematic development in the lifecycle of project (see Appendix B). The three major commits are: Starting data exploration with quality foundations, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Krunkal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>-</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Wallis</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> test implementation with hypothesis testing evidence and the final report integration making findings into a comprehensive documentation where the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>reslts</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> of the statistical tests can be reproduced and the findings communicated in a clear way.</w:t></w:r></w:p>'
$para2.InsertXML($xml2)

# --- Edit 3: Insert new "6. Conclusions" section before the trailing bookmark paragraph ---
$bm = $d.Bookmarks("_GoBack")
$bmPara = $bm.Range.Paragraphs(1).Range
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:bookmarkStart w:id="56" w:name="_Toc214915994"/><w:r><w:t>6. Conclusions</w:t></w:r><w:bookmarkEnd w:id="56"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:bookmarkStart w:id="57" w:name="_Toc214911788"/><w:bookmarkStart w:id="58" w:name="_Toc214915304"/><w:bookmarkStart w:id="59" w:name="_Toc214915995"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>6.1 Results explained</w:t></w:r><w:bookmarkEnd w:id="57"/><w:bookmarkEnd w:id="58"/><w:bookmarkEnd w:id="59"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Kruskal-</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Wallis</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> test provided statistically significant values (H = 14.27, p = 0.014), which means that the levels of evidence scores prove that they are distributed differently in terms of popular interest. Nevertheless, low correlation by Spearman (0.067) and boxplots characteristics indicate that this importance is conditioned by a very high interest rate in two conditions where the evidence is the most significant (score 5). The results of the evidence with a score of 0 to 4 demonstrate significantly overlapping distributions meaning that the discrepancies do not have a significant influence on the patterns of public attention within each of the categories.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:bookmarkStart w:id="60" w:name="_Toc214915996"/><w:bookmarkStart w:id="61" w:name="_Toc214911789"/><w:bookmarkStart w:id="62" w:name="_Toc214915305"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>6.2 Discussion</w:t></w:r><w:bookmarkEnd w:id="60"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkEnd w:id="61"/><w:bookmarkEnd w:id="62"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Results indicate that there is a threshold effect: it is only extraordinary scientific evidence that has a greater effect on the interest of the population. It implies that dramatic claims or media coverage instead of incremental improvements in evidence are the motivators of public attention, which is in line with the findings of Ishida et al. (2020). In the case of the cannabis health policy, the moderate-quality evidence is having trouble with sensationalized claims. Healthcare communicators will need to understand that evidence quality is not enough to engage the population; touching stories and smart media strategies are also important.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:bookmarkStart w:id="63" w:name="_Toc214911790"/><w:bookmarkStart w:id="64" w:name="_Toc214915306"/><w:bookmarkStart w:id="65" w:name="_Toc214915997"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">6.3 </w:t></w:r><w:bookmarkEnd w:id="63"/><w:bookmarkEnd w:id="64"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Limitations</w:t></w:r><w:bookmarkEnd w:id="65"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Limitations of the studies are the small sample size to obtain a high evidence score and the inability to measure the confounding factors such as media coverage and celebrity endorsement. The research of the future needs to explore the variations in time (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Sarvet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> et al., 2018), explore the aspects of attention which drive people, and the interventions to match the interest of people with the quality of evidence using stringent longitudinal studies.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="MS Mincho" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr></w:p>'
$bmPara.InsertXML($xml3)

Write-Host "Edits applied"
